# Actualización automática 2025-07-08 15:30:08
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" : row 10 ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M10").Value = -79.79000000000001
$ws1.Range("O10").Value = -1.77
$ws1.Range("P10").Value = -4.67

# --- Sheet "VENTA MENSUAL" : row 10 and row 32 ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F10").Value = -86.23
$ws2.Range("F32").Value = 237.19

# --- Sheet "CUMPLIMIENTO MENSUAL" : rows 9, 15, 17, 18 ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D9").Value = 22.86
$ws3.Range("E9").Value = 627.39
$ws3.Range("F9").Value = 0.03515570934256055

$ws3.Range("D15").Value = 119.09
$ws3.Range("E15").Value = 23339.73
$ws3.Range("F15").Value = 0.005076555427766615

$ws3.Range("D17").Value = 10.67
$ws3.Range("E17").Value = 1589.33
$ws3.Range("F17").Value = 0.00666875

$ws3.Range("D18").Value = 226.91
$ws3.Range("E18").Value = 33707.80607548726
$ws3.Range("F18").Value = 0.006686662693603864

$wb.Save()
